# Update column G ("K") values for rows 2-22 on Sheet1.
# This reflects regenerating save_data to use K (strikeouts) instead of Strike#,
# with recalculated std/mean and written s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 2
    16 = 2
    17 = 3
    18 = 1
    19 = 0
    20 = 3
    21 = 0
    22 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
